# Update India Insurance (Life) database rows with refreshed capital-structure figures.
# Rows 3 and 5 (SBI Life / HDFC Life) swap positions and rows 2-6 receive refreshed
# historical growth, margin, cash, debt and coverage metrics.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 ---
$ws.Range("D2").Value = 0.147
$ws.Range("E2").Value = 0.125
$ws.Range("F2").Value = 0.0756
$ws.Range("G2").Value = 0.03194954318363741
$ws.Range("H2").Value = 0.03194954318363741
$ws.Range("I2").Value = 0.03380516673339871
$ws.Range("J2").Value = 0.02812453553219106
$ws.Range("K2").Value = 579.1
$ws.Range("L2").Value = 0.02664072060135987
$ws.Range("M2").Value = 15.6
$ws.Range("N2").Value = 0.0003538682793381756
$ws.Range("O2").Value = 0.02693835261612847
$ws.Range("P2").Value = 15.6
$ws.Range("Q2").Value = 0.0003538682793381756
$ws.Range("R2").Value = 0.02693835261612847
$ws.Range("U2").Value = 387.2
$ws.Range("V2").Value = 0.008783192164086
$ws.Range("W2").Value = 0.1665668440972659
$ws.Range("X2").Value = 0.0755802358550296
$ws.Range("Y2").Value = 0.09098660824223626
$ws.Range("Z2").Value = 7.164851877854451
$ws.Range("AA2").Value = 0.2251038003229546
$ws.Range("AB2").Value = 0.07539144086193472
$ws.Range("AC2").Value = 0.1497240164994088
$ws.Range("AD2").Value = 128.9
$ws.Range("AE2").Value = 43.47784324709436
$ws.Range("AF2").Value = 172.3778432470943
$ws.Range("AG2").Value = -214.8221567529057
$ws.Range("AH2").Value = 0.00389496548643326
$ws.Range("AI2").Value = 0.04266168106009507
$ws.Range("AJ2").Value = -0.004896858978043923
$ws.Range("AK2").Value = -0.05880096884859121
$ws.Range("AL2").Value = 11.38
$ws.Range("AM2").Value = 11.38
$ws.Range("AN2").Value = 0.1647848720072817
$ws.Range("AO2").Value = 64.4639718804921
$ws.Range("AP2").Value = -0.2746271652820464
$ws.Range("AQ2").Value = 64.4639718804921

# --- Row 3 ---
$ws.Range("B3").Value = "SBI Life Insurance Company Limited (BSE:540719)"
$ws.Range("D3").Value = 0.229
$ws.Range("E3").Value = 0.151
$ws.Range("F3").Value = 0.116
$ws.Range("G3").Value = 0.03224118762670852
$ws.Range("H3").Value = 0.03224118762670852
$ws.Range("I3").Value = 0.03396987340376408
$ws.Range("J3").Value = 0.02882571364355038
$ws.Range("K3").Value = 219.1
$ws.Range("L3").Value = 0.02865738015826303
$ws.Range("M3").Value = -0
$ws.Range("N3").Value = -0
$ws.Range("O3").Value = -0
$ws.Range("P3").Value = -0
$ws.Range("Q3").Value = -0
$ws.Range("R3").Value = -0
$ws.Range("T3").ClearContents()
$ws.Range("U3").Value = 293.3
$ws.Range("V3").Value = 0.02368474179351556
$ws.Range("W3").Value = 0.1917892156862745
$ws.Range("X3").Value = 0.07555574871048525
$ws.Range("Y3").Value = 0.1162334669757892
$ws.Range("Z3").Value = 8.168996546750016
$ws.Range("AA3").Value = 0.2354771552117679
$ws.Range("AB3").Value = 0.07539252931112668
$ws.Range("AC3").Value = 0.1600846259006412
$ws.Range("AE3").Value = 41.81666445760872
$ws.Range("AF3").Value = 41.81666445760872
$ws.Range("AG3").Value = -251.4833355423913
$ws.Range("AH3").Value = 0.003365440542632167
$ws.Range("AI3").Value = 0.03083559516196468
$ws.Range("AJ3").Value = -0.02072889796460187
$ws.Range("AK3").Value = -0.2366196767066423
$ws.Range("AL3").Value = 1.15
$ws.Range("AM3").Value = 1.15
$ws.Range("AO3").Value = 225.5652173913043
$ws.Range("AP3").Value = -0.9195675571975694
$ws.Range("AQ3").Value = 225.5652173913043

# --- Row 4 ---
$ws.Range("G4").Value = 0.03042897998093422
$ws.Range("H4").Value = 0.03042897998093422
$ws.Range("I4").Value = 0.03412774070543375
$ws.Range("J4").Value = 0.02575402971922176
$ws.Range("K4").Value = 32.1
$ws.Range("L4").Value = 0.01224022878932316
$ws.Range("M4").Value = -0
$ws.Range("N4").Value = -0
$ws.Range("O4").Value = -0
$ws.Range("P4").Value = -0
$ws.Range("Q4").Value = -0
$ws.Range("R4").Value = -0
$ws.Range("T4").ClearContents()
$ws.Range("U4").Value = 0
$ws.Range("V4").Value = 0
$ws.Range("W4").Value = 0.1074656846334115
$ws.Range("X4").Value = 0.07608567994203019
$ws.Range("Y4").Value = 0.03138000469138126
$ws.Range("Z4").Value = 8.779241821663385
$ws.Range("AA4").Value = 0.2261008547873534
$ws.Range("AB4").Value = 0.07536921523434872
$ws.Range("AC4").Value = 0.1507316395530046
$ws.Range("AD4").Value = 47.3
$ws.Range("AE4").Value = 0
$ws.Range("AF4").Value = 47.3
$ws.Range("AG4").Value = 47.3
$ws.Range("AH4").Value = 0.01461319822046465
$ws.Range("AI4").Value = 0.1124851367419738
$ws.Range("AJ4").Value = 0.01461319822046465
$ws.Range("AK4").Value = 0.1124851367419738
$ws.Range("AL4").Value = 5.33
$ws.Range("AM4").Value = 5.33
$ws.Range("AN4").Value = 0.4145486415425066
$ws.Range("AO4").Value = 16.79174484052533
$ws.Range("AP4").Value = 0.4145486415425066
$ws.Range("AQ4").Value = 16.79174484052533

# --- Row 5 ---
$ws.Range("B5").Value = "HDFC Life Insurance Company Limited (BSE:540777)"
$ws.Range("D5").Value = 0.147
$ws.Range("E5").Value = 0.125
$ws.Range("F5").Value = 0.0756
$ws.Range("G5").Value = 0.03426054467239104
$ws.Range("H5").Value = 0.03426054467239104
$ws.Range("I5").Value = 0.03653906587694316
$ws.Range("J5").Value = 0.03170659189321804
$ws.Range("K5").Value = 182.4
$ws.Range("L5").Value = 0.03126124736490308
$ws.Range("M5").Value = -0
$ws.Range("N5").Value = -0
$ws.Range("O5").Value = -0
$ws.Range("P5").Value = -0
$ws.Range("Q5").Value = -0
$ws.Range("R5").Value = -0
$ws.Range("T5").ClearContents()
$ws.Range("U5").Value = 46.2
$ws.Range("V5").Value = 0.002469003847798204
$ws.Range("W5").Value = 0.2030276046304542
$ws.Range("X5").Value = 0.07560472299957396
$ws.Range("Y5").Value = 0.1274228816308802
$ws.Range("Z5").Value = 7.068143640707459
$ws.Range("AA5").Value = 0.2241067458585558
$ws.Range("AB5").Value = 0.07539035241274275
$ws.Range("AC5").Value = 0.148716393445813
$ws.Range("AD5").Value = 81.59999999999999
$ws.Range("AE5").Value = 1.392561638998828
$ws.Range("AF5").Value = 82.99256163899882
$ws.Range("AG5").Value = 36.79256163899882
$ws.Range("AH5").Value = 0.004415674088022173
$ws.Range("AI5").Value = 0.07231253770650797
$ws.Range("AJ5").Value = 0.001962396325952121
$ws.Range("AK5").Value = 0.03340246037091001
$ws.Range("AL5").Value = 2.32
$ws.Range("AM5").Value = 2.32
$ws.Range("AN5").Value = 0.3785260677357554
$ws.Range("AO5").Value = 91.5948275862069
$ws.Range("AP5").Value = 0.1706733294011718
$ws.Range("AQ5").Value = 91.5948275862069

# --- Row 6 ---
$ws.Range("D6").Value = 0.106
$ws.Range("E6").Value = -0.08529999999999999
$ws.Range("F6").Value = -0.083
$ws.Range("G6").Value = 0.0298684934424193
$ws.Range("H6").Value = 0.0298684934424193
$ws.Range("I6").Value = 0.03060061344346685
$ws.Range("J6").Value = 0.0262213737103912
$ws.Range("K6").Value = 145.5
$ws.Range("L6").Value = 0.02582213782455144
$ws.Range("M6").Value = 15.6
$ws.Range("N6").Value = 0.001591966691158462
$ws.Range("O6").Value = 0.1072164948453608
$ws.Range("P6").Value = 15.6
$ws.Range("Q6").Value = 0.001591966691158462
$ws.Range("R6").Value = 0.1072164948453608
$ws.Range("U6").Value = 47.7
$ws.Range("V6").Value = 0.004867744305657604
$ws.Range("W6").Value = 0.1413444725082572
$ws.Range("X6").Value = 0.07540077346463257
$ws.Range("Y6").Value = 0.06594369904362464
$ws.Range("Z6").Value = 5.78648757082424
$ws.Range("AA6").Value = 0.1517296530651162
$ws.Range("AB6").Value = 0.07539944829622185
$ws.Range("AC6").Value = 0.07633020476889431
$ws.Range("AE6").Value = 0.268617150486815
$ws.Range("AF6").Value = 0.268617150486815
$ws.Range("AG6").Value = -47.43138284951318
$ws.Range("AH6").Value = 0.00002741139963616967
$ws.Range("AI6").Value = 0.0002406384506020759
$ws.Range("AJ6").Value = -0.004863874924810599
$ws.Range("AK6").Value = -0.0443877745315006
$ws.Range("AL6").Value = 2.58
$ws.Range("AM6").Value = 2.58
$ws.Range("AO6").Value = 66.74418604651162
$ws.Range("AP6").Value = -0.2648628976569736
$ws.Range("AQ6").Value = 66.74418604651162
